$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2901895501337322
$ws.Range("C2").Value = 0.04252466767447061
$ws.Range("E2").Value = 0.4576337759361593
$ws.Range("F2").Value = 2.33141574223427
$ws.Range("G2").Value = 0.002459547653490467
$ws.Range("J2").Value = 0.06334602002779111
$ws.Range("K2").Value = 0.2558571640329035
$ws.Range("M2").Value = 0.3865600791276904
$ws.Range("O2").Value = 3.189339725869459
$ws.Range("B3").Value = 0.2588489828966374
$ws.Range("C3").Value = 0.03959513240245371
$ws.Range("E3").Value = 0.4480533624814242
$ws.Range("F3").Value = 2.322064375550255
$ws.Range("G3").Value = 0.002461919933249564
$ws.Range("J3").Value = 0.06387607611698343
$ws.Range("K3").Value = 0.2242407072877626
$ws.Range("M3").Value = 0.3658086672436909
$ws.Range("O3").Value = 3.217667543967266
$ws.Range("B4").Value = 0.2396227967419975
$ws.Range("C4").Value = 0.03778527064849868
$ws.Range("E4").Value = 0.4423786807834631
$ws.Range("F4").Value = 2.31759659313029
$ws.Range("G4").Value = 0.002463453323096666
$ws.Range("J4").Value = 0.06421977872631057
$ws.Range("K4").Value = 0.2047932749426735
$ws.Range("M4").Value = 0.3532120404919326
$ws.Range("O4").Value = 3.236867881036929
$ws.Range("B5").Value = 0.2317926763755338
$ws.Range("C5").Value = 0.03704498807490353
$ws.Range("E5").Value = 0.4401185156980816
$ws.Range("F5").Value = 2.316096277690761
$ws.Range("G5").Value = 0.002464097561305009
$ws.Range("J5").Value = 0.06436443402015257
$ws.Range("K5").Value = 0.1968599530749344
$ws.Range("M5").Value = 0.3481154488998399
$ws.Range("O5").Value = 3.245146269155157
$ws.Range("B6").Value = 0.2304927875227918
$ws.Range("C6").Value = 0.03692190002658435
$ws.Range("E6").Value = 0.4397463789561087
$ws.Range("F6").Value = 2.315866498993245
$ws.Range("G6").Value = 0.002464205708208681
$ws.Range("J6").Value = 0.06438873153370395
$ws.Range("K6").Value = 0.1955421398806578
$ws.Range("M6").Value = 0.34727138295694
$ws.Range("O6").Value = 3.246548307358438
$ws.Range("B7").Value = 0.2395171773854656
$ws.Range("C7").Value = 0.03777529800986912
$ws.Range("E7").Value = 0.4423479874861798
$ws.Range("F7").Value = 2.317575062320245
$ws.Range("G7").Value = 0.002463461932942235
$ws.Range("J7").Value = 0.06422171099144336
$ws.Range("K7").Value = 0.2046863165859065
$ws.Range("M7").Value = 0.3531431574459347
$ws.Range("O7").Value = 3.236977688033306
$ws.Range("B8").Value = 0.2793800253742234
$ws.Range("C8").Value = 0.04151689738684183
$ws.Range("E8").Value = 0.4542873775857004
$ws.Range("F8").Value = 2.327926955347252
$ws.Range("G8").Value = 0.002460349710202625
$ws.Range("J8").Value = 0.06352499978990611
$ws.Range("K8").Value = 0.2449632714407102
$ws.Range("M8").Value = 0.3793750580384696
$ws.Range("O8").Value = 3.198731960483997
$ws.Range("B9").Value = 0.3576717454994309
$ws.Range("C9").Value = 0.04876430401358789
$ws.Range("E9").Value = 0.4793470151554686
$ws.Range("F9").Value = 2.358340538330339
$ws.Range("G9").Value = 0.002454853399063452
$ws.Range("J9").Value = 0.06230326151378041
$ws.Range("K9").Value = 0.3236557640414333
$ws.Range("M9").Value = 0.4319581672656483
$ws.Range("O9").Value = 3.138082626731034
$ws.Range("B10").Value = 0.4152526815371118
$ws.Range("C10").Value = 0.05403247067032169
$ws.Range("E10").Value = 0.4987623214730945
$ws.Range("F10").Value = 2.386864053102059
$ws.Range("G10").Value = 0.002451181454875121
$ws.Range("J10").Value = 0.0614933654202634
$ws.Range("K10").Value = 0.3812804555095113
$ws.Range("M10").Value = 0.4712827911671127
$ws.Range("O10").Value = 3.102291214740319
$ws.Range("B11").Value = 0.4414582643679807
$ws.Range("C11").Value = 0.05641649907366286
$ws.Range("E11").Value = 0.5078130715178446
$ws.Range("F11").Value = 2.401185025792685
$ws.Range("G11").Value = 0.00244958972189274
$ws.Range("J11").Value = 0.06114388826795736
$ws.Range("K11").Value = 0.4074515163631247
$ws.Range("M11").Value = 0.4893220988905185
$ws.Range("O11").Value = 3.087916854141326
$ws.Range("B12").Value = 0.4513829838763286
$ws.Range("C12").Value = 0.05731743477400641
$ws.Range("E12").Value = 0.511271757288128
$ws.Range("F12").Value = 2.40680162713717
$ws.Range("G12").Value = 0.002448998225053085
$ws.Range("J12").Value = 0.06101426965359114
$ws.Range("K12").Value = 0.4173553404021106
$ws.Range("M12").Value = 0.4961745785251424
$ws.Range("O12").Value = 3.082748282058532
$ws.Range("B13").Value = 0.4492454702479449
$ws.Range("C13").Value = 0.0571234849359854
$ws.Range("E13").Value = 0.5105254737373031
$ws.Range("F13").Value = 2.405583381627622
$ws.Range("G13").Value = 0.002449125114497286
$ws.Range("J13").Value = 0.06104206439585891
$ws.Range("K13").Value = 0.4152226748408054
$ws.Range("M13").Value = 0.4946978270175322
$ws.Range("O13").Value = 3.083849203189487
$ws.Range("B14").Value = 0.4422747548662187
$ws.Range("C14").Value = 0.05649065683235222
$ws.Range("E14").Value = 0.5080969915446616
$ws.Range("F14").Value = 2.401643227274775
$ws.Range("G14").Value = 0.002449540833735524
$ws.Range("J14").Value = 0.06113316996042517
$ws.Range("K14").Value = 0.4082664445346325
$ws.Range("M14").Value = 0.4898854292540094
$ws.Range("O14").Value = 3.087486125110246
$ws.Range("B15").Value = 0.438005137742465
$ws.Range("C15").Value = 0.05610278975582617
$ws.Range("E15").Value = 0.5066135582543154
$ws.Range("F15").Value = 2.399254978568806
$ws.Range("G15").Value = 0.002449796938379749
$ws.Range("J15").Value = 0.06118932896053941
$ws.Range("K15").Value = 0.4040046813033769
$ws.Range("M15").Value = 0.4869404747576382
$ws.Range("O15").Value = 3.089749629359716
$ws.Range("B16").Value = 0.4135402823274319
$ws.Range("C16").Value = 0.05387641338053584
$ws.Range("E16").Value = 0.4981752265956629
$ws.Range("F16").Value = 2.385955225331742
$ws.Range("G16").Value = 0.00245128705662075
$ws.Range("J16").Value = 0.06151658552636441
$ws.Range("K16").Value = 0.3795692174986982
$ws.Range("M16").Value = 0.470106885229626
$ws.Range("O16").Value = 3.103269028824684
$ws.Range("B17").Value = 0.3985345668887987
$ws.Range("C17").Value = 0.05250737070420541
$ws.Range("E17").Value = 0.493054525898323
$ws.Range("F17").Value = 2.378140925521606
$ws.Range("G17").Value = 0.00245222130324857
$ws.Range("J17").Value = 0.06172219716287852
$ws.Range("K17").Value = 0.3645675954205672
$ws.Range("M17").Value = 0.459818362495362
$ws.Range("O17").Value = 3.112051550828795
$ws.Range("B18").Value = 0.3899048068745969
$ws.Range("C18").Value = 0.05171876109749007
$ws.Range("E18").Value = 0.4901298172589676
$ws.Range("F18").Value = 2.373772985294835
$ws.Range("G18").Value = 0.002452766063599526
$ws.Range("J18").Value = 0.06184224330641008
$ws.Range("K18").Value = 0.3559350666639602
$ws.Range("M18").Value = 0.4539148408584026
$ws.Range("O18").Value = 3.1172825275479
$ws.Range("B19").Value = 0.3869831264120194
$ws.Range("C19").Value = 0.051451551548098
$ws.Range("E19").Value = 0.4891430986201399
$ws.Range("F19").Value = 2.372315823583122
$ws.Range("G19").Value = 0.002452951783765486
$ws.Range("J19").Value = 0.06188319535531939
$ws.Range("K19").Value = 0.3530115687436819
$ws.Range("M19").Value = 0.4519184497595035
$ws.Range("O19").Value = 3.119084463474024
$ws.Range("B20").Value = 0.4001318371813909
$ws.Range("C20").Value = 0.05265322933297512
$ws.Range("E20").Value = 0.4935975035184867
$ws.Range("F20").Value = 2.378959663893042
$ws.Range("G20").Value = 0.002452121085100864
$ws.Range("J20").Value = 0.06170012486859644
$ws.Range("K20").Value = 0.3661649601679358
$ws.Range("M20").Value = 0.4609121290229794
$ws.Range("O20").Value = 3.111098055540026
$ws.Range("B21").Value = 0.4443221930487198
$ws.Range("C21").Value = 0.05667658421936039
$ws.Range("E21").Value = 0.5088094450450171
$ws.Range("F21").Value = 2.402795292315417
$ws.Range("G21").Value = 0.002449418421917605
$ws.Range("J21").Value = 0.06110633623749706
$ws.Range("K21").Value = 0.4103098403743672
$ws.Range("M21").Value = 0.4912983685953876
$ws.Range("O21").Value = 3.086410414772985
$ws.Range("B22").Value = 0.4732101296685016
$ws.Range("C22").Value = 0.05929530810712436
$ws.Range("E22").Value = 0.5189341019215306
$ws.Range("F22").Value = 2.419501505696871
$ws.Range("G22").Value = 0.002447717670521768
$ws.Range("J22").Value = 0.06073411814471896
$ws.Range("K22").Value = 0.4391224611307791
$ws.Range("M22").Value = 0.5112820749801941
$ws.Range("O22").Value = 3.071876888651985
$ws.Range("B23").Value = 0.4577916039974923
$ws.Range("C23").Value = 0.05789864734572348
$ws.Range("E23").Value = 0.5135136830434703
$ws.Range("F23").Value = 2.410481823238939
$ws.Range("G23").Value = 0.00244861940881333
$ws.Range("J23").Value = 0.06093132831717618
$ws.Range("K23").Value = 0.4237483035377068
$ws.Range("M23").Value = 0.5006050760224809
$ws.Range("O23").Value = 3.079487051506788
$ws.Range("B24").Value = 0.3994097198902864
$ws.Range("C24").Value = 0.05258729140554408
$ws.Range("E24").Value = 0.4933519633679921
$ws.Range("F24").Value = 2.378589124141456
$ws.Range("G24").Value = 0.00245216636971295
$ws.Range("J24").Value = 0.06171009802944116
$ws.Range("K24").Value = 0.3654428161290753
$ws.Range("M24").Value = 0.4604176013783601
$ws.Range("O24").Value = 3.111528564523269
$ws.Range("B25").Value = 0.3364802467327479
$ws.Range("C25").Value = 0.04681348835407562
$ws.Range("E25").Value = 0.4723914256497537
$ws.Range("F25").Value = 2.349028773876128
$ws.Range("G25").Value = 0.002456275723666835
$ws.Range("J25").Value = 0.06261834312722492
$ws.Range("K25").Value = 0.3023998745417771
$ws.Range("M25").Value = 0.417611221343364
$ws.Range("O25").Value = 3.152951578564554
